$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.8243289364228475
$ws.Range("E2").Value = 0.8243289364228475

# Row 3
$ws.Range("D3").Value = 0.00001517562940939917
$ws.Range("E3").Value = 0.00001517562940939917

# Row 4
$ws.Range("D4").Value = 0.007720695492901638
$ws.Range("E4").Value = 0.007720695492901638

# Row 5
$ws.Range("D5").Value = 0.000349016615380066
$ws.Range("E5").Value = 0.000349016615380066

# Row 6
$ws.Range("D6").Value = 0.8032215568833873
$ws.Range("E6").Value = 0.8032215568833873

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.3051696791851207
$ws.Range("E7").Value = 0.6948303208148794

# Row 8
$ws.Range("D8").Value = 0.9999987551882091
$ws.Range("E8").Value = 0.000001244811790868106

# Row 9
$ws.Range("D9").Value = 0.9909720097608554
$ws.Range("E9").Value = 0.009027990239144579

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.0001598077461993702
$ws.Range("E11").Value = 0.9998401922538006
$ws.Range("F11").Value = 1.331043124198914
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.8986412885892574
$ws.Range("E12").Value = 0.8986412885892574

# Row 13
$ws.Range("D13").Value = 0.002470371116379446
$ws.Range("E13").Value = 0.002470371116379446

# Row 14
$ws.Range("D14").Value = 0.0005515971648240394
$ws.Range("E14").Value = 0.0005515971648240394

# Row 15
$ws.Range("D15").Value = 0.0006575401089120316
$ws.Range("E15").Value = 0.0006575401089120316

# Row 16
$ws.Range("D16").Value = 0.7864184865984901
$ws.Range("E16").Value = 0.7864184865984901

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.06677099487235535
$ws.Range("E17").Value = 0.9332290051276446

# Row 18
$ws.Range("D18").Value = 0.9999999886498433
$ws.Range("E18").Value = 0.00000001135015670339357

# Row 19
$ws.Range("D19").Value = 0.995520671794768
$ws.Range("E19").Value = 0.004479328205232003

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.2216503757703994
$ws.Range("E21").Value = 0.7783496242296006
$ws.Range("F21").Value = 0.8054138422012329
$ws.Range("G21").Value = 0.6
